$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C9").Value = 0.5
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 3
$ws.Range("C10").Value = 0.5
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 3
$ws.Range("C11").Value = 0.5
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 3
$ws.Range("C12").Value = 0.5
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 3
$ws.Range("C13").Value = 0.5
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 3
$ws.Range("C14").Value = 0.5
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 3
$ws.Range("C15").Value = 0.5
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 3
$ws.Range("C16").Value = 0.5
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 3
$ws.Range("C17").Value = 0.5
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 3
$ws.Range("C18").Value = 0.5
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 3
$ws.Range("C19").Value = 0.5
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 3
$ws.Range("C20").Value = 0.5
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 3
$ws.Range("C21").Value = 0.5
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 3
$ws.Range("C22").Value = 0.5
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 3
$ws.Range("C23").Value = 0.5
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 3
$ws.Range("C24").Value = 0.5
$ws.Range("D24").Value = 1
$ws.Range("E24").Value = 3
$ws.Range("C25").Value = 0.5
$ws.Range("D25").Value = 1.5
$ws.Range("E25").Value = 3
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = 9
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = 5
$ws.Range("C28").Value = 0.5
$ws.Range("D28").Value = 1.5
$ws.Range("E28").Value = 4
$ws.Range("C29").Value = 0.1
$ws.Range("D29").Value = 0.3
$ws.Range("E29").Value = 1.5
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = 2
$ws.Range("E30").Value = 6
$ws.Range("C31").Value = 1
$ws.Range("D31").Value = 2
$ws.Range("E31").Value = 6
$ws.Range("C32").Value = 1
$ws.Range("D32").Value = 2
$ws.Range("E32").Value = 6
$ws.Range("C33").Value = 0.5
$ws.Range("D33").Value = 1
$ws.Range("E33").Value = 3
$ws.Range("C34").Value = 1
$ws.Range("D34").Value = 2
$ws.Range("E34").Value = 4
$ws.Range("C35").Value = 0.2
$ws.Range("D35").Value = 1
$ws.Range("E35").Value = 1.5
$ws.Range("C36").Value = 1
$ws.Range("D36").Value = 1.5
$ws.Range("E36").Value = 3
$ws.Range("C37").Value = 0.2
$ws.Range("D37").Value = 0.5
$ws.Range("E37").Value = 2
$ws.Range("C38").Value = 0.1
$ws.Range("D38").Value = 0.5
$ws.Range("E38").Value = 1
$ws.Range("C39").Value = 0.3
$ws.Range("D39").Value = 1
$ws.Range("E39").Value = 1.5
$ws.Range("C40").Value = 0.2
$ws.Range("D40").Value = 0.5
$ws.Range("E40").Value = 1
$ws.Range("C41").Value = 2
$ws.Range("D41").Value = 3
$ws.Range("E41").Value = 5
$ws.Range("C42").Value = 1
$ws.Range("D42").Value = 2.5
$ws.Range("E42").Value = 4
$ws.Range("C43").Value = 0.5
$ws.Range("D43").Value = 1
$ws.Range("E43").Value = 2.5
$ws.Range("C44").Value = 0.5
$ws.Range("D44").Value = 1
$ws.Range("E44").Value = 2.5
$ws.Range("C45").Value = 0.5
$ws.Range("D45").Value = 1
$ws.Range("E45").Value = 2.5
$ws.Range("C46").Value = 0.2
$ws.Range("D46").Value = 0.5
$ws.Range("E46").Value = 1
$ws.Range("C47").Value = 1
$ws.Range("D47").Value = 2
$ws.Range("E47").Value = 3
$ws.Range("C48").Value = 0.2
$ws.Range("D48").Value = 0.5
$ws.Range("E48").Value = 1
$ws.Range("C49").Value = 1
$ws.Range("D49").Value = 1.5
$ws.Range("E49").Value = 2.5
$ws.Range("C50").Value = 0.2
$ws.Range("D50").Value = 0.5
$ws.Range("E50").Value = 1
$ws.Range("C51").Value = 3
$ws.Range("D51").Value = 4
$ws.Range("E51").Value = 6
$ws.Range("C52").Value = 1
$ws.Range("D52").Value = 2
$ws.Range("E52").Value = 3
$ws.Range("C53").Value = 1
$ws.Range("D53").Value = 1.5
$ws.Range("E53").Value = 2
$ws.Range("C54").Value = 1
$ws.Range("D54").Value = 1.5
$ws.Range("E54").Value = 2
$ws.Range("C55").Value = 1
$ws.Range("D55").Value = 1.5
$ws.Range("E55").Value = 2
$ws.Range("C56").Value = 1
$ws.Range("D56").Value = 1.5
$ws.Range("E56").Value = 2
$ws.Range("C57").Value = 2
$ws.Range("D57").Value = 3
$ws.Range("E57").Value = 6
$ws.Range("C58").Value = 2
$ws.Range("D58").Value = 3
$ws.Range("E58").Value = 6
$ws.Range("C59").Value = 2
$ws.Range("D59").Value = 3
$ws.Range("E59").Value = 6

$ws.Range("J9").Value = "Database"
$ws.Range("J10").Value = "Database"
$ws.Range("J11").Value = "Database"
$ws.Range("J12").Value = "Database"
$ws.Range("J13").Value = "Database"
$ws.Range("J14").Value = "Database"
$ws.Range("J15").Value = "Database"
$ws.Range("J16").Value = "Database"
$ws.Range("J17").Value = "Database"
$ws.Range("J18").Value = "Database"
$ws.Range("J19").Value = "Database"
$ws.Range("J20").Value = "Database"
$ws.Range("J21").Value = "Database"
$ws.Range("J22").Value = "Database"
$ws.Range("J23").Value = "Database"
$ws.Range("J32").Value = "Database"
$ws.Range("J33").Value = "Database"
$ws.Range("J34").Value = "Database"
$ws.Range("J35").Value = "Database"
$ws.Range("J36").Value = "Database"

$ws.Range("D65").Select()
